$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab ("Waypoints" -> "WayPoints")
$ws.Name = "WayPoints"

# Fix the DEVAC row longitude formatting
$ws.Range("D3").Value = "W087°26'06.06"""

# Replace the old HERKO/ALABAMA waypoint row with the new ERLIN/GEORGIA waypoint
$ws.Range("B5").Value = "ERLIN"
$ws.Range("C5").Value = "N34°05'13.95"""
$ws.Range("D5").Value = "W085°01'18.94"""
$ws.Range("F5").Value = "GEORGIA"

# Add the new KATL/26L destination row
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "KATL/26L"
$ws.Range("C6").Value = "obtained from runways data"
$ws.Range("D6").Value = "obtained from runways data"
$ws.Range("E6").Value = "USA"
$ws.Range("F6").Value = "ATLANTA"

# Column width tweaks (waypoint name column widened, new bearing column sized)
$ws.Columns.Item(3).ColumnWidth = 23.98307291666667
$ws.Columns.Item(6).ColumnWidth = 8.346354166666666

# Restore the active-cell selection left after editing
$ws.Range("H5").Select()
